# Applies the Thu Dec 28 15:31:36 UTC 2023 cryptos-list price/volume update
# (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is a cell and its new text value. "Numeric" marks Price-column
# values that look like plain numbers; for those we force a text NumberFormat
# before assigning so Excel keeps them as text (preserving things like the
# trailing zero in "1.00"/"1.50") instead of silently converting to a Double,
# then we restore the original (default) style so no visible formatting changes.
$updates = @(
    @{ Cell = 'D2'; Value = '42.623.14'; Numeric = $false }
    @{ Cell = 'E2'; Value = '  -0.89%  '; Numeric = $false }
    @{ Cell = 'D3'; Value = '2.359.93'; Numeric = $false }
    @{ Cell = 'E3'; Value = '  +1.01%  '; Numeric = $false }
    @{ Cell = 'E4'; Value = '  -0.06%  '; Numeric = $false }
    @{ Cell = 'D5'; Value = '327.11'; Numeric = $true }
    @{ Cell = 'E5'; Value = '  +5.73%  '; Numeric = $false }
    @{ Cell = 'D6'; Value = '99.76'; Numeric = $true }
    @{ Cell = 'E6'; Value = '  -7.43%  '; Numeric = $false }
    @{ Cell = 'D7'; Value = '0.635'; Numeric = $true }
    @{ Cell = 'E7'; Value = '  +0.17%  '; Numeric = $false }
    @{ Cell = 'E8'; Value = '  +0.16%  '; Numeric = $false }
    @{ Cell = 'D9'; Value = '0.615'; Numeric = $true }
    @{ Cell = 'E9'; Value = '  -1.36%  '; Numeric = $false }
    @{ Cell = 'D10'; Value = '39.99'; Numeric = $true }
    @{ Cell = 'E10'; Value = '  -7.59%  '; Numeric = $false }
    @{ Cell = 'D11'; Value = '0.0918'; Numeric = $true }
    @{ Cell = 'E11'; Value = '  -1.98%  '; Numeric = $false }
    @{ Cell = 'D12'; Value = '8.35'; Numeric = $true }
    @{ Cell = 'E12'; Value = '  -6.40%  '; Numeric = $false }
    @{ Cell = 'D13'; Value = '1.01'; Numeric = $true }
    @{ Cell = 'E13'; Value = '  -5.96%  '; Numeric = $false }
    @{ Cell = 'E14'; Value = '  +0.53%  '; Numeric = $false }
    @{ Cell = 'D15'; Value = '16.19'; Numeric = $true }
    @{ Cell = 'E15'; Value = '  -0.34%  '; Numeric = $false }
    @{ Cell = 'D16'; Value = '2.727.22'; Numeric = $false }
    @{ Cell = 'E16'; Value = '  +1.19%  '; Numeric = $false }
    @{ Cell = 'D17'; Value = '2.366.83'; Numeric = $false }
    @{ Cell = 'E17'; Value = '  +0.92%  '; Numeric = $false }
    @{ Cell = 'D18'; Value = '42.637.79'; Numeric = $false }
    @{ Cell = 'E18'; Value = '  -0.72%  '; Numeric = $false }
    @{ Cell = 'D19'; Value = '7.68'; Numeric = $true }
    @{ Cell = 'E19'; Value = '  +6.37%  '; Numeric = $false }
    @{ Cell = 'E20'; Value = '  -2.17%  '; Numeric = $false }
    @{ Cell = 'B21'; Value = 'Litecoin'; Numeric = $false }
    @{ Cell = 'C21'; Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; Numeric = $false }
    @{ Cell = 'D21'; Value = '75.03'; Numeric = $true }
    @{ Cell = 'E21'; Value = '  -0.47%  '; Numeric = $false }
    @{ Cell = 'B22'; Value = 'PancakeSwap'; Numeric = $false }
    @{ Cell = 'C22'; Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; Numeric = $false }
    @{ Cell = 'D22'; Value = '3.69'; Numeric = $true }
    @{ Cell = 'E22'; Value = '  +6.87%  '; Numeric = $false }
    @{ Cell = 'D23'; Value = '275.25'; Numeric = $true }
    @{ Cell = 'E23'; Value = '  +10.15%  '; Numeric = $false }
    @{ Cell = 'D24'; Value = '2.31'; Numeric = $true }
    @{ Cell = 'E24'; Value = '  -8.34%  '; Numeric = $false }
    @{ Cell = 'D25'; Value = '9.61'; Numeric = $true }
    @{ Cell = 'E25'; Value = '  +7.59%  '; Numeric = $false }
    @{ Cell = 'D26'; Value = '1.00'; Numeric = $true }
    @{ Cell = 'E26'; Value = '  +0.03%  '; Numeric = $false }
    @{ Cell = 'D27'; Value = '11.40'; Numeric = $true }
    @{ Cell = 'E27'; Value = '  -4.00%  '; Numeric = $false }
    @{ Cell = 'D28'; Value = '23.84'; Numeric = $true }
    @{ Cell = 'E28'; Value = '  +6.57%  '; Numeric = $false }
    @{ Cell = 'E29'; Value = '  -1.64%  '; Numeric = $false }
    @{ Cell = 'D30'; Value = '173.68'; Numeric = $true }
    @{ Cell = 'E30'; Value = '  -0.18%  '; Numeric = $false }
    @{ Cell = 'E31'; Value = '  -1.61%  '; Numeric = $false }
    @{ Cell = 'D32'; Value = '0.0898'; Numeric = $true }
    @{ Cell = 'E32'; Value = '  -0.93%  '; Numeric = $false }
    @{ Cell = 'D33'; Value = '35.09'; Numeric = $true }
    @{ Cell = 'E33'; Value = '  -9.22%  '; Numeric = $false }
    @{ Cell = 'E34'; Value = '  +0.40%  '; Numeric = $false }
    @{ Cell = 'E35'; Value = '  +0.83%  '; Numeric = $false }
    @{ Cell = 'D36'; Value = '4.58'; Numeric = $true }
    @{ Cell = 'E36'; Value = '  -8.85%  '; Numeric = $false }
    @{ Cell = 'D37'; Value = '0.0357'; Numeric = $true }
    @{ Cell = 'E37'; Value = '  -4.88%  '; Numeric = $false }
    @{ Cell = 'D38'; Value = '3.89'; Numeric = $true }
    @{ Cell = 'E38'; Value = '  -5.40%  '; Numeric = $false }
    @{ Cell = 'D39'; Value = '2.87'; Numeric = $true }
    @{ Cell = 'E39'; Value = '  +4.09%  '; Numeric = $false }
    @{ Cell = 'E40'; Value = '  +0.45%  '; Numeric = $false }
    @{ Cell = 'D41'; Value = '1.50'; Numeric = $true }
    @{ Cell = 'E41'; Value = '  +1.44%  '; Numeric = $false }
    @{ Cell = 'D42'; Value = '0.226'; Numeric = $true }
    @{ Cell = 'E42'; Value = '  -2.70%  '; Numeric = $false }
    @{ Cell = 'D43'; Value = '68.30'; Numeric = $true }
    @{ Cell = 'E43'; Value = '  -4.90%  '; Numeric = $false }
    @{ Cell = 'E44'; Value = '  -0.04%  '; Numeric = $false }
    @{ Cell = 'D45'; Value = '89.37'; Numeric = $true }
    @{ Cell = 'E45'; Value = '  +40.72%  '; Numeric = $false }
    @{ Cell = 'D46'; Value = '114.96'; Numeric = $true }
    @{ Cell = 'E46'; Value = '  +4.97%  '; Numeric = $false }
    @{ Cell = 'D47'; Value = '11.88'; Numeric = $true }
    @{ Cell = 'E47'; Value = '  -4.01%  '; Numeric = $false }
    @{ Cell = 'D48'; Value = '5.45'; Numeric = $true }
    @{ Cell = 'E48'; Value = '  -4.05%  '; Numeric = $false }
    @{ Cell = 'D49'; Value = '8.97'; Numeric = $true }
    @{ Cell = 'E49'; Value = '  -2.81%  '; Numeric = $false }
    @{ Cell = 'D50'; Value = '1.596.07'; Numeric = $false }
    @{ Cell = 'E50'; Value = '  +7.23%  '; Numeric = $false }
    @{ Cell = 'E51'; Value = '  -2.75%  '; Numeric = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Numeric) {
        $origStyle = $rng.Style
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = $origStyle
    } else {
        $rng.Value = $u.Value
    }
}
